# Update of results using corrected definition of APR.
# Only the "exp_arms" sheet holds hard-coded numbers; the "consort" sheet
# pulls everything via formulas (e.g. =exp_arms!C16) so it recalculates
# automatically once these source values change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp_arms")

# Row 16
$ws.Range("B16").Value = 1770
$ws.Range("C16").Value = 1954
$ws.Range("D16").Value = 1470
$ws.Range("E16").Value = 2580
$ws.Range("F16").Value = 1904

# Row 17
$ws.Range("B17").Value = 1386
$ws.Range("C17").Value = 1467
$ws.Range("D17").Value = 1176
$ws.Range("E17").Value = 1982
$ws.Range("F17").Value = 1534
$ws.Range("G17").Value = 6

# Row 18
$ws.Range("B18").Value = 1386
$ws.Range("C18").Value = 1467
$ws.Range("D18").Value = 1176
$ws.Range("E18").Value = 1982
$ws.Range("F18").Value = 1534

# Row 23
$ws.Range("B23").Value = 19246
$ws.Range("C23").Value = 19348
$ws.Range("E23").Value = 19583
